$wb = $excel.ActiveWorkbook

# sheet1 (Worksheets.Item(1))
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 62
$ws.Range("F4").Value = 5124
$ws.Range("F5").Value = 5124
$ws.Range("F6").Value = 121
$ws.Range("F7").Value = 157
$ws.Range("F8").Value = 207
$ws.Range("F9").Value = 43
$ws.Range("F11").Value = 174
$ws.Range("F12").Value = 8502
$ws.Range("F13").Value = 8502
$ws.Range("F16").Value = 619
$ws.Range("F17").Value = 2555
$ws.Range("F18").Value = 6330
$ws.Range("F20").Value = 8
$ws.Range("F22").Value = 2529
$ws.Range("F23").Value = 22
$ws.Range("F24").Value = 15
$ws.Range("F25").Value = 6449
$ws.Range("F26").Value = 198
$ws.Range("F27").Value = 72
$ws.Range("F28").Value = 139
$ws.Range("F30").Value = 459
$ws.Range("F31").Value = 6956
$ws.Range("F33").Value = 34
$ws.Range("F34").Value = 232
$ws.Range("F35").Value = 13
$ws.Range("F36").Value = 12
$ws.Range("F38").Value = 8
$ws.Range("F40").Value = 2
$ws.Range("F42").Value = 50
$ws.Range("F43").Value = 2536
$ws.Range("F45").Value = 70
$ws.Range("F46").Value = 1130
$ws.Range("F47").Value = 64
$ws.Range("F48").Value = 530
$ws.Range("F49").Value = 2454
$ws.Range("F50").Value = 83
$ws.Range("F51").Value = 1124

# sheet2 (Worksheets.Item(2))
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 183
$ws.Range("F6").Value = 75
$ws.Range("F7").Value = 18

# sheet4 (Worksheets.Item(4))
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 62
$ws.Range("F3").Value = 5124
$ws.Range("F4").Value = 5124
$ws.Range("F5").Value = 121
$ws.Range("F6").Value = 157
$ws.Range("F7").Value = 207
$ws.Range("F8").Value = 43
$ws.Range("F10").Value = 174
$ws.Range("F11").Value = 8502
$ws.Range("F12").Value = 8502
$ws.Range("F15").Value = 619
$ws.Range("F16").Value = 2555
$ws.Range("F17").Value = 183
$ws.Range("F18").Value = 6330
$ws.Range("F20").Value = 75
$ws.Range("F21").Value = 8
$ws.Range("F22").Value = 2529
$ws.Range("F23").Value = 22
$ws.Range("F24").Value = 18
$ws.Range("F26").Value = 15
$ws.Range("F27").Value = 6449
$ws.Range("F28").Value = 198
$ws.Range("F29").Value = 72
$ws.Range("F30").Value = 139
$ws.Range("F32").Value = 459
$ws.Range("F33").Value = 6956
$ws.Range("F35").Value = 34
$ws.Range("F36").Value = 232
$ws.Range("F37").Value = 12
$ws.Range("F41").Value = 50
$ws.Range("F42").Value = 2536
$ws.Range("F44").Value = 70
$ws.Range("F45").Value = 1130
$ws.Range("F46").Value = 64
$ws.Range("F47").Value = 530
$ws.Range("F49").Value = 2455
$ws.Range("F50").Value = 83
$ws.Range("F51").Value = 1124
